$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$ws.Cells.Item(3, 2).Value = "6.0.0"

# Date update
$ws.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty -> "Alvearie Team"
$ws.Cells.Item(9, 2).Value = "Alvearie Team"

# Row 10 was duplicate "Contact" -> now "Jurisdiction" / "United States of America"
$ws.Cells.Item(10, 1).Value = "Jurisdiction"
$ws.Cells.Item(10, 2).Value = "United States of America"

# Row 11 was duplicate "Contact" -> now "Description" / full description text
$ws.Cells.Item(11, 1).Value = "Description"
$ws.Cells.Item(11, 2).Value = "The Value Set for the flag indicating the type of room and board services on a facility Claim"

# Row 12 was "Description" / text -> now "Purpose" / empty
$ws.Cells.Item(12, 1).Value = "Purpose"
$ws.Cells.Item(12, 2).Value = ""

# Row 13 was "Purpose" / empty -> now "Copyright" / empty
$ws.Cells.Item(13, 1).Value = "Copyright"
$ws.Cells.Item(13, 2).Value = ""

# Row 14 was "Copyright" / empty -> now "Immutable" / "BooleanType[null]"
$ws.Cells.Item(14, 1).Value = "Immutable"
$ws.Cells.Item(14, 2).Value = "BooleanType[null]"

# Row 15 (was "Immutable" / "BooleanType[null]") no longer exists - delete entire row
$ws.Rows.Item(15).Delete()
